$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 737.25
$ws.Range("I9").Value = 154.85715
$ws.Range("K9").Value = 154.85715
$ws.Range("M9").Value = 14.14285000000001
$ws.Range("H19").Value = 1939.5
$ws.Range("I19").Value = 2739.8572
$ws.Range("J19").Value = 819
$ws.Range("K19").Value = 2739.8572
$ws.Range("L19").Value = 819
$ws.Range("M19").Value = -2564.8572
$ws.Range("N19").Value = -1169
$ws.Range("H80").Value = 3373.4736
$ws.Range("J80").Value = 3956.6
$ws.Range("L80").Value = 11869.8
$ws.Range("N80").Value = -13865.8
$ws.Range("H82").Value = 811.5
$ws.Range("I82").Value = 811.5
$ws.Range("K82").Value = 2434.5
$ws.Range("M82").Value = -2028.5
$ws.Range("H83").Value = 3373.4736
$ws.Range("J83").Value = 3956.6
$ws.Range("L83").Value = 35609.4
$ws.Range("N83").Value = -45593.4
$ws.Range("H85").Value = 811.5
$ws.Range("I85").Value = 811.5
$ws.Range("K85").Value = 2434.5
$ws.Range("M85").Value = -1030.5
$ws.Range("H92").Value = 773.6316
$ws.Range("I92").Value = 599.94116
$ws.Range("K92").Value = 599.94116
$ws.Range("M92").Value = 648.05884
$ws.Range("H98").Value = 2811.0967
$ws.Range("I98").Value = 978.61536
$ws.Range("J98").Value = 12340
$ws.Range("K98").Value = 978.61536
$ws.Range("L98").Value = 12340
$ws.Range("M98").Value = 519.38464
$ws.Range("N98").Value = -15336
$ws.Range("H103").Value = 2710
$ws.Range("J103").Value = 2710
$ws.Range("L103").Value = 8130
$ws.Range("N103").Value = -9302
$ws.Range("H117").Value = 100742
$ws.Range("J117").Value = 100742
$ws.Range("L117").Value = 100742
$ws.Range("N117").Value = -109920
$ws.Range("H118").Value = 843.875
$ws.Range("I118").Value = 907.2857
$ws.Range("K118").Value = 2721.8571
$ws.Range("M118").Value = -1064.8571
$ws.Range("H122").Value = 2811.0967
$ws.Range("I122").Value = 978.61536
$ws.Range("J122").Value = 12340
$ws.Range("K122").Value = 2935.84608
$ws.Range("L122").Value = 37020
$ws.Range("M122").Value = -485.8460800000003
$ws.Range("N122").Value = -41920
$ws.Range("H125").Value = 5466.647
$ws.Range("I125").Value = 4565.8
$ws.Range("J125").Value = 5842
$ws.Range("K125").Value = 41092.2
$ws.Range("L125").Value = 52578
$ws.Range("M125").Value = -38632.2
$ws.Range("N125").Value = -57498
$ws.Range("H132").Value = 5981.7607
$ws.Range("I132").Value = 5055.892
$ws.Range("J132").Value = 9788.111000000001
$ws.Range("K132").Value = 15167.676
$ws.Range("L132").Value = 29364.333
$ws.Range("M132").Value = -12637.676
$ws.Range("N132").Value = -34424.333
$ws.Range("H135").Value = 3402.742
$ws.Range("I135").Value = 1558.4642
$ws.Range("K135").Value = 14026.1778
$ws.Range("M135").Value = -11491.1778
$ws.Range("H138").Value = 224527.1
$ws.Range("I138").Value = 1194.8077
$ws.Range("J138").Value = 530139.7
$ws.Range("K138").Value = 3584.4231
$ws.Range("L138").Value = 1590419.1
$ws.Range("M138").Value = 1555.5769
$ws.Range("N138").Value = -1600699.1
$ws.Range("H141").Value = 3066.8262
$ws.Range("I141").Value = 1946.8334
$ws.Range("J141").Value = 7098.8
$ws.Range("K141").Value = 5840.5002
$ws.Range("L141").Value = 21296.4
$ws.Range("M141").Value = -660.5002000000004
$ws.Range("N141").Value = -31656.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1029.5
$ws.Range("I22").Value = 733.3333
$ws.Range("J22").Value = 1918
$ws.Range("K22").Value = 733.3333
$ws.Range("L22").Value = 1918
$ws.Range("M22").Value = -434.3333
$ws.Range("N22").Value = -2516
$ws.Range("H61").Value = 4813.604
$ws.Range("I61").Value = 3375.0244
$ws.Range("K61").Value = 3375.0244
$ws.Range("M61").Value = -3163.0244
$ws.Range("H63").Value = 2214.6667
$ws.Range("I63").Value = 2214.6667
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2214.6667
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1528.6667
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2214.6667
$ws.Range("I66").Value = 2214.6667
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 11073.3335
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -7641.333500000001
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 2051.7083
$ws.Range("I74").Value = 2091.2632
$ws.Range("K74").Value = 2091.2632
$ws.Range("M74").Value = -1217.2632
$ws.Range("H77").Value = 2051.7083
$ws.Range("I77").Value = 2091.2632
$ws.Range("K77").Value = 10456.316
$ws.Range("M77").Value = -6088.315999999999
$ws.Range("H88").Value = 2374.5
$ws.Range("J88").Value = 1811.75
$ws.Range("L88").Value = 1811.75
$ws.Range("N88").Value = -2623.75
$ws.Range("H91").Value = 2374.5
$ws.Range("J91").Value = 1811.75
$ws.Range("L91").Value = 1811.75
$ws.Range("N91").Value = -4619.75
$ws.Range("H102").Value = 2150.6775
$ws.Range("I102").Value = 2167.9656
$ws.Range("K102").Value = 2167.9656
$ws.Range("M102").Value = -545.9656
$ws.Range("H132").Value = 4826.3667
$ws.Range("I132").Value = 2011.625
$ws.Range("K132").Value = 6034.875
$ws.Range("M132").Value = -3504.875
$ws.Range("H136").Value = 4813.604
$ws.Range("I136").Value = 3375.0244
$ws.Range("K136").Value = 10125.0732
$ws.Range("M136").Value = -7575.073199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H36").Value = 1632.6666
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H39").Value = 10973
$ws.Range("J39").Value = 17898
$ws.Range("L39").Value = 17898
$ws.Range("N39").Value = -18676
$ws.Range("H44").Value = 15500
$ws.Range("J44").Value = 21000
$ws.Range("L44").Value = 21000
$ws.Range("N44").Value = -21994
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H99").Value = 2550.1667
$ws.Range("J99").Value = 2999.5
$ws.Range("L99").Value = 2999.5
$ws.Range("N99").Value = -5995.5
$ws.Range("H105").Value = 2397.577
$ws.Range("I105").Value = 1442.2858
$ws.Range("K105").Value = 1442.2858
$ws.Range("M105").Value = 304.7141999999999
$ws.Range("H107").Value = 1669.1666
$ws.Range("I107").Value = 1472.5143
$ws.Range("J107").Value = 2652.4285
$ws.Range("K107").Value = 1472.5143
$ws.Range("L107").Value = 2652.4285
$ws.Range("M107").Value = 447.4857
$ws.Range("N107").Value = -6492.4285
$ws.Range("H120").Value = 97690.25
$ws.Range("J120").Value = 97690.25
$ws.Range("L120").Value = 97690.25
$ws.Range("N120").Value = -107366.25
$ws.Range("H134").Value = 2783.7727
$ws.Range("I134").Value = 2285.6667
$ws.Range("K134").Value = 6857.000100000001
$ws.Range("M134").Value = -4322.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2870.85
$ws.Range("I31").Value = 1994.6923
$ws.Range("K31").Value = 1994.6923
$ws.Range("M31").Value = -1699.6923
$ws.Range("H34").Value = 2870.85
$ws.Range("I34").Value = 1994.6923
$ws.Range("K34").Value = 1994.6923
$ws.Range("M34").Value = -1792.6923
$ws.Range("H132").Value = 4356.206
$ws.Range("I132").Value = 3041.5925
$ws.Range("K132").Value = 9124.7775
$ws.Range("M132").Value = -6594.7775
$ws.Range("H134").Value = 1507.3684
$ws.Range("I134").Value = 1452.2222
$ws.Range("K134").Value = 4356.6666
$ws.Range("M134").Value = -1821.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2278.5881
$ws.Range("I14").Value = 2278.5881
$ws.Range("K14").Value = 6835.7643
$ws.Range("M14").Value = -6662.7643
$ws.Range("H34").Value = 205
$ws.Range("J34").Value = 1000
$ws.Range("L34").Value = 3000
$ws.Range("N34").Value = -3168
$ws.Range("H107").Value = 4475.3335
$ws.Range("J107").Value = 4712.857
$ws.Range("L107").Value = 14138.571
$ws.Range("N107").Value = -17978.571
$ws.Range("H116").Value = 3073.4
$ws.Range("J116").Value = 2863.125
$ws.Range("L116").Value = 8589.375
$ws.Range("N116").Value = -15473.375
$ws.Range("H117").Value = 1666
$ws.Range("I117").Value = 905.6
$ws.Range("K117").Value = 2716.8
$ws.Range("M117").Value = 725.1999999999998
$ws.Range("H120").Value = 15377.6
$ws.Range("I120").Value = 8962.666999999999
$ws.Range("J120").Value = 25000
$ws.Range("K120").Value = 26888.001
$ws.Range("L120").Value = 75000
$ws.Range("M120").Value = -22050.001
$ws.Range("N120").Value = -84676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 35794
$ws.Range("J26").Value = 38392
$ws.Range("L26").Value = 38392
$ws.Range("N26").Value = -38952
$ws.Range("H45").Value = 70162.5
$ws.Range("J45").Value = 70162.5
$ws.Range("L45").Value = 70162.5
$ws.Range("N45").Value = -71280.5
$ws.Range("H50").Value = 35794
$ws.Range("J50").Value = 38392
$ws.Range("L50").Value = 38392
$ws.Range("N50").Value = -39388
$ws.Range("H63").Value = 79990.60000000001
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 79990.60000000001
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 79990.60000000001
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -81362.60000000001
$ws.Range("H66").Value = 79990.60000000001
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 79990.60000000001
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 239971.8
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -246835.8
$ws.Range("H80").Value = 4833.1113
$ws.Range("I80").Value = 5087.25
$ws.Range("J80").Value = 2800
$ws.Range("K80").Value = 5087.25
$ws.Range("L80").Value = 2800
$ws.Range("M80").Value = -4089.25
$ws.Range("N80").Value = -4796
$ws.Range("H83").Value = 4833.1113
$ws.Range("I83").Value = 5087.25
$ws.Range("J83").Value = 2800
$ws.Range("K83").Value = 25436.25
$ws.Range("L83").Value = 14000
$ws.Range("M83").Value = -20444.25
$ws.Range("N83").Value = -23984
$ws.Range("H102").Value = 22703.72
$ws.Range("I102").Value = 2786.739
$ws.Range("J102").Value = 251749
$ws.Range("K102").Value = 2786.739
$ws.Range("L102").Value = 251749
$ws.Range("M102").Value = -1164.739
$ws.Range("N102").Value = -254993
$ws.Range("H107").Value = 697.8
$ws.Range("J107").Value = 806.3333
$ws.Range("L107").Value = 806.3333
$ws.Range("N107").Value = -4646.3333
$ws.Range("H113").Value = 5432.375
$ws.Range("I113").Value = 5744.8335
$ws.Range("K113").Value = 5744.8335
$ws.Range("M113").Value = -3574.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3663.3684
$ws.Range("J22").Value = 3874.8235
$ws.Range("L22").Value = 3874.8235
$ws.Range("N22").Value = -4464.8235
$ws.Range("H27").Value = 3663.3684
$ws.Range("J27").Value = 3874.8235
$ws.Range("L27").Value = 3874.8235
$ws.Range("N27").Value = -4088.8235
$ws.Range("H33").Value = 21228.5
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H50").Value = 33688
$ws.Range("I50").Value = 31033
$ws.Range("J50").Value = 38998
$ws.Range("K50").Value = 31033
$ws.Range("L50").Value = 38998
$ws.Range("M50").Value = -30396
$ws.Range("N50").Value = -40272
$ws.Range("H54").Value = 44143
$ws.Range("J54").Value = 44143
$ws.Range("L54").Value = 44143
$ws.Range("N54").Value = -45431

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 5371.6665
$ws.Range("J21").Value = 5646.2
$ws.Range("L21").Value = 5646.2
$ws.Range("N21").Value = -6116.2
$ws.Range("H35").Value = 5371.6665
$ws.Range("J35").Value = 5646.2
$ws.Range("L35").Value = 5646.2
$ws.Range("N35").Value = -6226.2
$ws.Range("H40").Value = 14282.875
$ws.Range("I40").Value = 9253.75
$ws.Range("J40").Value = 19312
$ws.Range("K40").Value = 9253.75
$ws.Range("L40").Value = 19312
$ws.Range("M40").Value = -9104.75
$ws.Range("N40").Value = -19610
$ws.Range("H48").Value = 36472.668
$ws.Range("H50").Value = 29728
$ws.Range("J50").Value = 29728
$ws.Range("L50").Value = 29728
$ws.Range("N50").Value = -30990
$ws.Range("H55").Value = 7758.8887
$ws.Range("I55").Value = 1641.5
$ws.Range("J55").Value = 19993.666
$ws.Range("K55").Value = 1641.5
$ws.Range("L55").Value = 19993.666
$ws.Range("M55").Value = -1364.5
$ws.Range("N55").Value = -20547.666
$ws.Range("H62").Value = 115850.29
$ws.Range("I62").Value = 134075.33
$ws.Range("K62").Value = 134075.33
$ws.Range("M62").Value = -133451.33
$ws.Range("H64").Value = 99993
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H65").Value = 115850.29
$ws.Range("I65").Value = 134075.33
$ws.Range("K65").Value = 670376.6499999999
$ws.Range("M65").Value = -667256.6499999999
$ws.Range("H67").Value = 99993
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H100").Value = 1214.15
$ws.Range("I100").Value = 1124.125
$ws.Range("J100").Value = 1574.25
$ws.Range("K100").Value = 2248.25
$ws.Range("L100").Value = 3148.5
$ws.Range("M100").Value = -1707.25
$ws.Range("N100").Value = -4230.5
$ws.Range("H122").Value = 2021.4878
$ws.Range("I122").Value = 1894.3226
$ws.Range("K122").Value = 5682.9678
$ws.Range("M122").Value = -3232.9678
$ws.Range("H126").Value = 4905.4165
$ws.Range("I126").Value = 4988.091
$ws.Range("J126").Value = 3996
$ws.Range("K126").Value = 14964.273
$ws.Range("L126").Value = 11988
$ws.Range("M126").Value = -12494.273
$ws.Range("N126").Value = -16928
$ws.Range("H132").Value = 1442.4186
$ws.Range("I132").Value = 1067.2941
$ws.Range("K132").Value = 3201.8823
$ws.Range("M132").Value = -671.8823000000002
$ws.Range("H133").Value = 52101.777
$ws.Range("J133").Value = 52101.777
$ws.Range("L133").Value = 52101.777
$ws.Range("N133").Value = -62221.777
